$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (scheme labels), rows 3-19, in final desired order ---
$ws.Range("B3").Value = "ND Single"
$ws.Range("B4").Value = "RD Single"
$ws.Range("B5").Value = "TD Single"
$ws.Range("B6").Value = "Morris"
$ws.Range("B7").Value = "Ring Perpendicular to ND"
$ws.Range("B8").Value = "Ring Perpendicular to RD"
$ws.Range("B9").Value = "Ring Perpendicular to TD"
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# --- Column A (index numbers) for new rows 17-19 (copy style from A16) ---
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A16").Copy($ws.Range("A18"))
$ws.Range("A16").Copy($ws.Range("A19"))
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# --- Numeric data grid, columns C:M, rows 3-19 ---
$ws.Range("C3").Value = 1.03
$ws.Range("D3").Value = 0.89
$ws.Range("E3").Value = 1.01
$ws.Range("F3").Value = 1.03
$ws.Range("G3").Value = 0.92
$ws.Range("H3").Value = 1.1
$ws.Range("I3").Value = 1.02
$ws.Range("J3").Value = 0.89
$ws.Range("K3").Value = 0.95
$ws.Range("L3").Value = 0.99
$ws.Range("M3").Value = 0.9949999999999998
$ws.Range("C4").Value = 1.01
$ws.Range("D4").Value = 0.95
$ws.Range("E4").Value = 1.02
$ws.Range("F4").Value = 1.01
$ws.Range("G4").Value = 0.97
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.95
$ws.Range("K4").Value = 0.985
$ws.Range("L4").Value = 0.9975
$ws.Range("M4").Value = 0.9916666666666667
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.99
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.99
$ws.Range("J5").Value = 0.99
$ws.Range("K5").Value = 0.995
$ws.Range("L5").Value = 0.9975
$ws.Range("M5").Value = 0.9966666666666667
$ws.Range("C6").Value = 0.98
$ws.Range("D6").Value = 1.04
$ws.Range("E6").Value = 0.99
$ws.Range("F6").Value = 0.98
$ws.Range("G6").Value = 1.01
$ws.Range("H6").Value = 0.98
$ws.Range("I6").Value = 0.99
$ws.Range("J6").Value = 1.04
$ws.Range("K6").Value = 1.015
$ws.Range("L6").Value = 0.9974999999999999
$ws.Range("M6").Value = 0.9983333333333334
$ws.Range("C7").Value = 1.01
$ws.Range("D7").Value = 0.953013698630137
$ws.Range("E7").Value = 1.011780821917808
$ws.Range("F7").Value = 1.01
$ws.Range("G7").Value = 0.9680821917808219
$ws.Range("H7").Value = 1.020821917808219
$ws.Range("I7").Value = 1.004657534246575
$ws.Range("J7").Value = 0.953013698630137
$ws.Range("K7").Value = 0.9823972602739726
$ws.Range("L7").Value = 0.9961986301369863
$ws.Range("M7").Value = 0.9947260273972601
$ws.Range("C8").Value = 0.9994736842105263
$ws.Range("D8").Value = 0.9757894736842105
$ws.Range("E8").Value = 0.9978947368421053
$ws.Range("F8").Value = 0.9994736842105263
$ws.Range("G8").Value = 0.9847368421052631
$ws.Range("H8").Value = 1.013157894736842
$ws.Range("I8").Value = 0.9994736842105263
$ws.Range("J8").Value = 0.9757894736842105
$ws.Range("K8").Value = 0.986842105263158
$ws.Range("L8").Value = 0.993157894736842
$ws.Range("M8").Value = 0.9950877192982457
$ws.Range("C9").Value = 0.9957894736842106
$ws.Range("D9").Value = 0.9947368421052631
$ws.Range("E9").Value = 0.998421052631579
$ws.Range("F9").Value = 0.9957894736842106
$ws.Range("G9").Value = 0.9936842105263158
$ws.Range("H9").Value = 1.001578947368421
$ws.Range("I9").Value = 0.9957894736842106
$ws.Range("J9").Value = 0.9947368421052631
$ws.Range("K9").Value = 0.996578947368421
$ws.Range("L9").Value = 0.9961842105263159
$ws.Range("M9").Value = 0.9966666666666667
$ws.Range("C10").Value = 1.002839692827024
$ws.Range("D10").Value = 0.9662521697578933
$ws.Range("E10").Value = 1.006223748755438
$ws.Range("F10").Value = 1.002839692827024
$ws.Range("G10").Value = 0.9781768555407804
$ws.Range("H10").Value = 1.015563277946242
$ws.Range("I10").Value = 1.000406761700286
$ws.Range("J10").Value = 0.9662521697578933
$ws.Range("K10").Value = 0.9862379592566659
$ws.Range("L10").Value = 0.9945388260418451
$ws.Range("M10").Value = 0.9949104177546108
$ws.Range("C11").Value = 0.992743152300807
$ws.Range("D11").Value = 1.004192933784116
$ws.Range("E11").Value = 0.9966661581427609
$ws.Range("F11").Value = 0.992743152300807
$ws.Range("G11").Value = 1.000106479251719
$ws.Range("H11").Value = 0.9937855693156546
$ws.Range("I11").Value = 0.9937293540626293
$ws.Range("J11").Value = 1.004192933784116
$ws.Range("K11").Value = 1.000429545963438
$ws.Range("L11").Value = 0.9965863491321226
$ws.Range("M11").Value = 0.9968706078096145
$ws.Range("C12").Value = 0.9926968014870151
$ws.Range("D12").Value = 1.004302551118491
$ws.Range("E12").Value = 0.9966482894307098
$ws.Range("F12").Value = 0.9926968014870151
$ws.Range("G12").Value = 1.000181284590684
$ws.Range("H12").Value = 0.9937705815575472
$ws.Range("I12").Value = 0.9937256434911336
$ws.Range("J12").Value = 1.004302551118491
$ws.Range("K12").Value = 1.0004754202746
$ws.Range("L12").Value = 0.9965861108808076
$ws.Range("M12").Value = 0.9968875252792634
$ws.Range("C13").Value = 0.9927328123083401
$ws.Range("D13").Value = 1.004179197928528
$ws.Range("E13").Value = 0.996705264006105
$ws.Range("F13").Value = 0.9927328123083401
$ws.Range("G13").Value = 1.000116944162057
$ws.Range("H13").Value = 0.993784006567703
$ws.Range("I13").Value = 0.9937411545587762
$ws.Range("J13").Value = 1.004179197928528
$ws.Range("K13").Value = 1.000442230967316
$ws.Range("L13").Value = 0.9965875216378282
$ws.Range("M13").Value = 0.9968765632552515
$ws.Range("C14").Value = 1.011679999999999
$ws.Range("D14").Value = 0.9462320000000016
$ws.Range("E14").Value = 1.00354
$ws.Range("F14").Value = 1.011679999999999
$ws.Range("G14").Value = 0.9585239999999999
$ws.Range("H14").Value = 1.050023999999998
$ws.Range("I14").Value = 1.007088000000001
$ws.Range("J14").Value = 0.9462320000000016
$ws.Range("K14").Value = 0.974886000000001
$ws.Range("L14").Value = 0.9932830000000001
$ws.Range("M14").Value = 0.9961813333333334
$ws.Range("C15").Value = 1.03
$ws.Range("D15").Value = 0.89
$ws.Range("E15").Value = 1.01
$ws.Range("F15").Value = 1.03
$ws.Range("G15").Value = 0.92
$ws.Range("H15").Value = 1.1
$ws.Range("I15").Value = 1.02
$ws.Range("J15").Value = 0.89
$ws.Range("K15").Value = 0.95
$ws.Range("L15").Value = 0.99
$ws.Range("M15").Value = 0.9949999999999998
$ws.Range("C16").Value = 1.015982313267198
$ws.Range("D16").Value = 0.9335898371072014
$ws.Range("E16").Value = 1.005099884543999
$ws.Range("F16").Value = 1.015982313267198
$ws.Range("G16").Value = 0.9524768059392041
$ws.Range("H16").Value = 1.055838604287997
$ws.Range("I16").Value = 1.010153071616001
$ws.Range("J16").Value = 0.9335898371072014
$ws.Range("K16").Value = 0.9693448608256001
$ws.Range("L16").Value = 0.9926635870463989
$ws.Range("M16").Value = 0.9955234194602666
$ws.Range("C17").Value = 0.9957241291045527
$ws.Range("D17").Value = 0.9962966581934911
$ws.Range("E17").Value = 0.9962310770160026
$ws.Range("F17").Value = 0.9957241291045527
$ws.Range("G17").Value = 0.9959585781772423
$ws.Range("H17").Value = 0.9964776362511404
$ws.Range("I17").Value = 0.9961099636175197
$ws.Range("J17").Value = 0.9962966581934911
$ws.Range("K17").Value = 0.9962638676047468
$ws.Range("L17").Value = 0.9959939983546497
$ws.Range("M17").Value = 0.9961330070599915
$ws.Range("C18").Value = 0.9957292296509429
$ws.Range("D18").Value = 1.002299830798997
$ws.Range("E18").Value = 0.9945436101798945
$ws.Range("F18").Value = 0.9957292296509429
$ws.Range("G18").Value = 0.9990045352957896
$ws.Range("H18").Value = 0.993653374055565
$ws.Range("I18").Value = 0.994527121314956
$ws.Range("J18").Value = 1.002299830798997
$ws.Range("K18").Value = 0.9984217204894459
$ws.Range("L18").Value = 0.9970754750701943
$ws.Range("M18").Value = 0.9966262835493575
$ws.Range("C19").Value = 0.993433849132324
$ws.Range("D19").Value = 1.009890087684527
$ws.Range("E19").Value = 0.9928980448015077
$ws.Range("F19").Value = 0.993433849132324
$ws.Range("G19").Value = 1.003858497415755
$ws.Range("H19").Value = 0.9866742492194698
$ws.Range("I19").Value = 0.9926562014478908
$ws.Range("J19").Value = 1.009890087684527
$ws.Range("K19").Value = 1.001394066243017
$ws.Range("L19").Value = 0.9974139576876706
$ws.Range("M19").Value = 0.9965684882835791